# Applies the Xhosa translation edits described by the diff.
# Each paragraph's entire run text is matched verbatim via Find, then the
# found Range's .Text is set directly (rather than using Find's Replace
# parameter) so that straight apostrophes in the new text are not
# auto-corrected into curly "smart quotes" by the host.
#
# NOTE: expressions are assigned to intermediate variables before being
# passed into Replace-ExactText; calling the function with two inline
# parenthesized expressions (e.g. `Foo (expr1) (expr2)`) is not reliable
# in this interpreter, so we avoid that pattern entirely.

$d = $word.ActiveDocument

$RSQUO = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK
$LSQUO = [char]0x2018   # U+2018 LEFT SINGLE QUOTATION MARK

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    $rng.Text = $newText
}

$old = "Appendix 16: SWIFT Interview Guide: Engagement  "
$new = "ISihlomelo 16: ISikhokelo Sodliwano-ndlebe lwe-SWIFT: Uthethathethwano  "
Replace-ExactText $old $new

$old = "Briefing:"
$new = "Ingxelo:"
Replace-ExactText $old $new

$old = "Hi there. Thank you for making the time for this phone call; it won" + $RSQUO + "t take longer than 15 minutes. We noticed that you didn" + $RSQUO + "t finish working through the ParentText programme. We would like to hear a little bit more about your experience, so we can hopefully improve the programme. "
$new = "Molo apho. Thank you for making the time for this phone call; it won" + $RSQUO + "t take longer than 15 minutes. Siqaphele ukuba akhange ugqibe ukusebenza kwinkqubo ye-ParentText. Singathanda ukuva ngakumbi kancinci malunga namava akho, ukuze sibenethemba lokuphucula inkqubo. "
Replace-ExactText $old $new

$old = "There are no right or wrong answers. You can skip any questions you do not feel comfortable answering. You can also stop this conversation at any time if you wish. If you decide at a later stage that you would like your contribution to be removed from the study, you can contact the research team by email until the [date to be determined]. "
$new = "Akukho zimpendulo zilungileyo okanye zingalunganga. Ungatsiba nayiphi na imibuzo oziva ungakhululekanga ukuyiphendula. Ungakwazi nokuyimisa lencoko nanini na ukuba uyafuna. Ukuba uthatha isgqibo sekumva sokuba ungathanda igalelo lakho lisuswe koluphononongo, ungaqhagamshelana neqela lophando nge-imeyile kude kube [ngumhla oza kumiselwa]. "
Replace-ExactText $old $new

$old = "Do you understand what I" + $RSQUO + "ve just explained? Do you have any questions? Can we begin?"
$new = "Uyayiqonda lento ndigqiba ukuyichaza? Ingaba unayo nayiphi na imibuzo? Singaqalisa?"
Replace-ExactText $old $new

$old = "What are the factors that you think led to you not finishing the programme?"
$new = "Zeziphi izinto ocinga ukuba zikukhokelele ekubeni ungayigqibi inkqubo?"
Replace-ExactText $old $new

$old = "What was your experience of the content?"
$new = "Ebenjani amava akho ngomxholo?"
Replace-ExactText $old $new

$old = "Probe: Explore relevance of content"
$new = "Buza: Phonononga ukufikeleleka komxholo"
Replace-ExactText $old $new

$old = "Probe: Explore interest in content"
$new = "Buza: Phonononga umdla kumxholo"
Replace-ExactText $old $new

$old = "Probe: What content could we have added to improve your experience? "
$new = "Buza: Ngowuphi umxholo engesiwengezile ukuphucula amava akho? "
Replace-ExactText $old $new

$old = "What was your experience of the time and data needed to complete the programme?"
$new = "Ebenjani amava akho ngexesha kunye nedatha efunekayo ukuze ugqibe inkqubo?"
Replace-ExactText $old $new

$old = "Probe for data: Explore the use of Wi-Fi hotspots in the community and what they used to connect to ParentText; explore customisation of content delivery e.g. whether they used audio/visual only."
$new = "Buza nge datha: Phonononga ukusetyenziswa kwe-Wi-Fi hotspots ekuhlaleni kunye nento ababeyisebenzisela ukunxulumelana ne-ParentText; phonononga ukulungiswa kokuhanjiswa komxholo umz. nokuba basebenzise iaudio/imiboniso kuphela."
Replace-ExactText $old $new

$old = "What was your experience of the home exercises/activities, and do you feel like they impacted whether you finished the programme?"
$new = "Ebesithini amava akho kwimisebenzi yasekhaya/ imisebenzi, kwaye uziva ngathi ibe nefuthe ekubeni uyigqibile inkqubo?"
Replace-ExactText $old $new

$old = "What else can you recommend we do to improve the programme?"
$new = "Yintoni engenye ongayicebisa siyenze ukuphucula lenkqubo?"
Replace-ExactText $old $new

$old = "Debriefing"
$new = "Ingxoxo"
Replace-ExactText $old $new

$old = "Is there anything that we haven" + $RSQUO + "t spoken about that you" + $RSQUO + "d like us to know? If after this conversation, there is anything that causes you to worry, remember you can still access the referral services in the programme by typing " + $LSQUO + "help" + $RSQUO + ". SADAG might be a good option. If you have any other questions about the study, you can send us a WhatsApp or email. "
$new = "Ingaba ikhona nantoni na esingakhange sithethe malunga nayo ongathanda siyazi? Ukuba emva kwalencoko, kukho nantoni na ekwenza ukuba ukhathazeke, khumbula usengafikelela kwiinkonzo zokuthunyelwa kwinkqubo ngokubhala 'nceda'. South African Depression and Anxiety Group (SADAG) inokuba lukhetho olulungileyo. If you have any other questions about the study, you can send us a WhatsApp or email. "
Replace-ExactText $old $new

$old = "We" + $RSQUO + "ve come to the end of our interview"
$new = "Sifikelele esiphelweni sodliwano-ndlebe lwethu"
Replace-ExactText $old $new

$old = "Thanks again for taking the time to speak with us today. Your responses have been very helpful! "
$new = "Enkosi kwakhona ngokuthatha ixesha lokuthetha nathi namhlanje. Iimpendulo zakho zibe luncedo kakhulu! "
Replace-ExactText $old $new

Write-Output "Applied all Xhosa translation replacements."
